$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.69"
$ws.Range("E2").Value = "'0.99%"
$ws.Range("G2").Value = "'21"
$ws.Range("D3").Value = "'27.22"
$ws.Range("E3").Value = "'0.93%"
$ws.Range("G3").Value = "'21"
$ws.Range("D4").Value = "'4.708"
$ws.Range("E4").Value = "'0.56%"
$ws.Range("G4").Value = "'21"
$ws.Range("D5").Value = "'0.06210"
$ws.Range("E5").Value = "'3.23%"
$ws.Range("G5").Value = "'21"
$ws.Range("D6").Value = "'6.719"
$ws.Range("E6").Value = "'0.61%"
$ws.Range("G6").Value = "'21"
$ws.Range("D7").Value = "'0.8511"
$ws.Range("E7").Value = "'-0.83%"
$ws.Range("G7").Value = "'21"
$ws.Range("D8").Value = "'0.9145"
$ws.Range("E8").Value = "'-1.31%"
$ws.Range("G8").Value = "'21"
$ws.Range("D9").Value = "'0.1409"
$ws.Range("E9").Value = "'0.83%"
$ws.Range("G9").Value = "'21"
$ws.Range("D10").Value = "'0.04626"
$ws.Range("E10").Value = "'-9.91%"
$ws.Range("G10").Value = "'21"
$ws.Range("D11").Value = "'0.07082"
$ws.Range("E11").Value = "'0.76%"
$ws.Range("G11").Value = "'21"
$ws.Range("D12").Value = "'0.03135"
$ws.Range("E12").Value = "'0.40%"
$ws.Range("G12").Value = "'21"
$ws.Range("D13").Value = "'0.09055"
$ws.Range("E13").Value = "'-0.79%"
$ws.Range("G13").Value = "'21"
$ws.Range("D14").Value = "'0.001546"
$ws.Range("E14").Value = "'1.11%"
$ws.Range("G14").Value = "'21"
$ws.Range("D15").Value = "'0.0006141"
$ws.Range("E15").Value = "'1.65%"
$ws.Range("G15").Value = "'21"
$ws.Range("D16").Value = "'0.005997"
$ws.Range("E16").Value = "'-0.83%"
$ws.Range("G16").Value = "'21"
$ws.Range("E17").Value = "'0.11%"
$ws.Range("G17").Value = "'21"
$ws.Range("D18").Value = "'3.171"
$ws.Range("E18").Value = "'-0.01%"
$ws.Range("G18").Value = "'21"
$ws.Range("E19").Value = "'-0.29%"
$ws.Range("G19").Value = "'21"
$ws.Range("E20").Value = "'0.43%"
$ws.Range("G20").Value = "'21"
$ws.Range("D21").Value = "'0.1308"
$ws.Range("E21").Value = "'0.80%"
$ws.Range("G21").Value = "'21"
$ws.Range("E22").Value = "'-1.04%"
$ws.Range("G22").Value = "'21"
$ws.Range("D23").Value = "'0.04228"
$ws.Range("E23").Value = "'0.04%"
$ws.Range("G23").Value = "'21"
$ws.Range("D24").Value = "'0.001209"
$ws.Range("E24").Value = "'-0.74%"
$ws.Range("G24").Value = "'21"
$ws.Range("E25").Value = "'-5.85%"
$ws.Range("G25").Value = "'21"
$ws.Range("D26").Value = "'0.0001200"
$ws.Range("E26").Value = "'-0.03%"
$ws.Range("G26").Value = "'21"
$ws.Range("E27").Value = "'5.04%"
$ws.Range("G27").Value = "'21"
$ws.Range("G28").Value = "'21"
$ws.Range("G29").Value = "'21"
$ws.Range("G30").Value = "'21"
$ws.Range("G31").Value = "'21"
$ws.Range("G32").Value = "'21"
$ws.Range("G33").Value = "'21"
$ws.Range("G34").Value = "'21"
$ws.Range("G35").Value = "'21"
$ws.Range("G36").Value = "'21"
$ws.Range("G37").Value = "'21"
$ws.Range("G38").Value = "'21"
$ws.Range("G39").Value = "'21"
$ws.Range("D40").Value = "'0.03910"
$ws.Range("E40").Value = "'1.73%"
$ws.Range("G40").Value = "'21"
$ws.Range("D41").Value = "'0.1114"
$ws.Range("E41").Value = "'-0.19%"
$ws.Range("G41").Value = "'21"
$ws.Range("D42").Value = "'0.004132"
$ws.Range("E42").Value = "'3.95%"
$ws.Range("G42").Value = "'21"
$ws.Range("D43").Value = "'0.002161"
$ws.Range("E43").Value = "'-1.81%"
$ws.Range("G43").Value = "'21"
$ws.Range("D44").Value = "'0.01392"
$ws.Range("E44").Value = "'-8.86%"
$ws.Range("G44").Value = "'21"
$ws.Range("D45").Value = "'0.00005172"
$ws.Range("E45").Value = "'1.86%"
$ws.Range("G45").Value = "'21"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("G46").Value = "'21"
$ws.Range("D47").Value = "'0.03590"
$ws.Range("E47").Value = "'-34.22%"
$ws.Range("G47").Value = "'21"
$ws.Range("D48").Value = "'0.1667"
$ws.Range("E48").Value = "'23.17%"
$ws.Range("G48").Value = "'21"
$ws.Range("E49").Value = "'0.02%"
$ws.Range("G49").Value = "'21"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("G50").Value = "'21"
$ws.Range("G51").Value = "'21"
